$d = $word.ActiveDocument

# Move to the very end of the document body.
$r = $d.Content
$r.Collapse(0)

$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$paraCountBefore = $d.Paragraphs.Count

# Build the whole block to append in a single InsertXML call (the runtime's
# Range object does not auto-extend across repeated InsertXML calls, so we
# batch everything together instead of re-inserting at a stale position).
$xml = ""

# 1) Page break paragraph (mirrors the default pPrDefault spacing that Word
#    stamps on a freshly inserted paragraph) followed by one blank paragraph.
$xml += "<w:p $w><w:pPr><w:spacing w:after=`"160`" w:line=`"259`" w:lineRule=`"auto`"/></w:pPr><w:r><w:br w:type=`"page`"/></w:r></w:p>"
$xml += "<w:p $w/>"

# 2) Heading "Nombre de test vs class" (inserted plain text; the Heading 1
#    style gets applied afterwards via the Paragraph object, since InsertXML
#    silently drops a bare <w:pStyle> in this runtime).
$xml += "<w:p $w><w:r><w:t>Nombre de test vs class</w:t></w:r></w:p>"

# 3) The French (fr-CA) statistic lines.
$xml += "<w:p $w><w:pPr><w:rPr><w:lang w:val=`"fr-CA`"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=`"fr-CA`"/></w:rPr><w:t>Nombre de fichier test:359</w:t></w:r></w:p>"
$xml += "<w:p $w><w:pPr><w:rPr><w:lang w:val=`"fr-CA`"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=`"fr-CA`"/></w:rPr><w:t>Nombre de public class:436</w:t></w:r></w:p>"
$xml += "<w:p $w><w:pPr><w:rPr><w:lang w:val=`"fr-CA`"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=`"fr-CA`"/></w:rPr><w:t xml:space=`"preserve`">Nombre de </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:rPr><w:lang w:val=`"fr-CA`"/></w:rPr><w:t>private</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:rPr><w:lang w:val=`"fr-CA`"/></w:rPr><w:t xml:space=`"preserve`"> class:6</w:t></w:r></w:p>"
$xml += "<w:p $w><w:pPr><w:rPr><w:lang w:val=`"fr-CA`"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=`"fr-CA`"/></w:rPr><w:t>Nombre de interface:114</w:t></w:r></w:p>"
$xml += "<w:p $w><w:pPr><w:rPr><w:lang w:val=`"fr-CA`"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=`"fr-CA`"/></w:rPr><w:t>Nombre de abstract class:33</w:t></w:r></w:p>"
$xml += "<w:p $w><w:pPr><w:rPr><w:lang w:val=`"fr-CA`"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=`"fr-CA`"/></w:rPr><w:t xml:space=`"preserve`">Nombre de autre( </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:rPr><w:lang w:val=`"fr-CA`"/></w:rPr><w:t>example</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:rPr><w:lang w:val=`"fr-CA`"/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:rPr><w:lang w:val=`"fr-CA`"/></w:rPr><w:t>enum</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:rPr><w:lang w:val=`"fr-CA`"/></w:rPr><w:t>) :75</w:t></w:r></w:p>"

$r.InsertXML($xml)

# Re-apply the Heading 1 paragraph style to the heading paragraph we added
# (it's the 3rd new paragraph: break-paragraph, blank paragraph, heading).
$headingIndex = $paraCountBefore + 3
$headingPara = $d.Paragraphs.Item($headingIndex)
$headingPara.Style = "Heading 1"

Write-Output "done"
